# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Refresh the COVID numbers for a handful of countries (Rusia, Filipinas,
#   Israel, Polonia, Singapur, Croacia, Hungria, Hong Kong, Eslovaquia,
#   Estonia)
# - Re-sort the table by "Casos totales" (column B) descending, which is how
#   this sheet is always kept ordered; Hungria's update pushes it above
#   Haiti and Finlandia.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados..." banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 10:45"

function Set-CountryRow {
    param(
        [string]$Country,
        $Total,
        $Nuevos,
        $Activos,
        $Recuperados,
        $Criticos,
        $Muertes
    )

    $finder = $ws.Range("A4:A219")
    $cell = $finder.Find($Country)
    if ($cell -eq $null) {
        Write-Host "WARNING: country not found: $Country"
        return
    }
    $r = $cell.Row

    if ($Total -ne $null)       { $ws.Cells.Item($r, 2).Value = $Total }
    if ($Nuevos -ne $null)      { $ws.Cells.Item($r, 3).Value = $Nuevos }
    if ($Activos -ne $null)     { $ws.Cells.Item($r, 4).Value = $Activos }
    if ($Recuperados -ne $null) { $ws.Cells.Item($r, 5).Value = $Recuperados }
    if ($Criticos -ne $null)    { $ws.Cells.Item($r, 7).Value = $Criticos }
    if ($Muertes -ne $null)     { $ws.Cells.Item($r, 8).Value = $Muertes }
}

Set-CountryRow "Rusia"      1025505 5195 840949 166736 61  17820
Set-CountryRow "Filipinas"  237365  2839 184687 48803  85  3875
Set-CountryRow "Israel"     129349  413  102107 26232  3   1010
Set-CountryRow "Polonia"    70824   437  54256  14448  7   2120
Set-CountryRow "Singapur"   57022   40   $null  728    $null $null
Set-CountryRow "Croacia"    11964   225  9008   2758   1   198
Set-CountryRow "Hungria"    8387    495  3958   3805   $null 624
Set-CountryRow "Hong Kong"  $null   $null 4493  271    $null $null
Set-CountryRow "Eslovaquia" 4614    88   2802   1775   $null $null
Set-CountryRow "Estonia"    2516    25   2170   282    $null $null

# The data rows (4..219) are always kept sorted by "Casos totales" (column B)
# descending; re-apply that after the updates above so any country whose
# total jumped past its neighbours (Hungria) moves to its new rank.
$sortRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4:B219")
$sortRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 2)
